$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status på uppgifterna var inte uppdaterad - sätt till "Klar" istället för "Ej Färdig"
$ws.Range("D2").Value = "Klar"
$ws.Range("D3").Value = "Klar"

# Uppdatera markerad cell till D3
$ws.Range("D3").Select()
